$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 (di) values from 0 to the new measured values
$ws.Range("B5").Value = 350
$ws.Range("C5").Value = 276
$ws.Range("D5").Value = 500
$ws.Range("E5").Value = 430
$ws.Range("F5").Value = 387

# Add new row 7: "hi" method values
$ws.Range("A7").Value = "hi"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 4

# Add new row 8: "ri" method values
$ws.Range("A8").Value = "ri"
$ws.Range("B8").Value = 23
$ws.Range("C8").Value = 45
$ws.Range("D8").Value = 34
$ws.Range("E8").Value = 38
$ws.Range("F8").Value = 12

# Update the active selection to match the author's final cursor position
$ws.Range("M8").Select()
